$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.441.21'
$ws.Range("E2").Value = '  -1.55%  '
$ws.Range("D3").Value = '3.008.95'
$ws.Range("E3").Value = '  -1.88%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '585.38'
$ws.Range("E5").Value = '  -0.60%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.86'
$ws.Range("E6").Value = '  -5.15%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.520'
$ws.Range("E8").Value = '  -3.40%  '
$ws.Range("D9").Value = '3.006.86'
$ws.Range("E9").Value = '  -1.88%  '
$ws.Range("E10").Value = '  -4.00%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.76'
$ws.Range("E11").Value = '  -0.98%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.441'
$ws.Range("E12").Value = '  -2.18%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000229'
$ws.Range("E13").Value = '  -3.37%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.74'
$ws.Range("E14").Value = '  -5.99%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.122'
$ws.Range("E15").Value = '  +2.43%  '
$ws.Range("D16").Value = '3.509.65'
$ws.Range("E16").Value = '  -1.82%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '62.452.18'
$ws.Range("E17").Value = '  -1.63%  '
$ws.Range("B18").Value = 'Polkadot'
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.01'
$ws.Range("E18").Value = '  -1.57%  '
$ws.Range("D19").Value = '3.011.45'
$ws.Range("E19").Value = '  -1.74%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '464.51'
$ws.Range("E20").Value = '  -1.71%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.93'
$ws.Range("E21").Value = '  -2.91%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.684'
$ws.Range("E22").Value = '  -2.99%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.34'
$ws.Range("E23").Value = '  -2.44%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.33'
$ws.Range("E24").Value = '  -3.98%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '79.94'
$ws.Range("E25").Value = '  -0.93%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.35'
$ws.Range("E26").Value = '  -3.62%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.19'
$ws.Range("E27").Value = '  -2.18%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  +0.01%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("E30").Value = '  -1.29%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.14'
$ws.Range("E31").Value = '  -4.68%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.12'
$ws.Range("E32").Value = '  -1.14%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.47'
$ws.Range("E33").Value = '  +1.24%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.106'
$ws.Range("E34").Value = '  -4.83%  '
$ws.Range("E35").Value = '  -0.79%  '
$ws.Range("D36").Value = '0.0₃0793'
$ws.Range("E36").Value = '  -3.36%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.73'
$ws.Range("E37").Value = '  -4.34%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.13'
$ws.Range("E38").Value = '  -3.58%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '50.34'
$ws.Range("E39").Value = '  -0.41%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.94'
$ws.Range("E40").Value = '  -3.13%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.93'
$ws.Range("E41").Value = '  -11.95%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '420.62'
$ws.Range("E42").Value = '  -3.89%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.112'
$ws.Range("E43").Value = '  +0.86%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.277'
$ws.Range("E44").Value = '  -3.47%  '
$ws.Range("D45").Value = '2.781.28'
$ws.Range("E45").Value = '  -0.42%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0353'
$ws.Range("E46").Value = '  -1.66%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '37.73'
$ws.Range("E47").Value = '  -7.91%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '129.37'
$ws.Range("E48").Value = '  -0.55%  '
$ws.Range("E49").Value = '  +0.01%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '24.07'
$ws.Range("E50").Value = '  -4.23%  '
$ws.Range("E51").Value = '  -1.19%  '
